$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date in column C for rows 2-16 from 45208 to 45212
for ($row = 2; $row -le 16; $row++) {
    $ws.Cells.Item($row, 3).Value = 45212
}

# Update the hyperlink formulas in row 2 to include the extra filename suffixes
$ws.Range("S2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_1885/artfynd/A 34293-2023 artfynd.xlsx", "A 34293-2023")'
$ws.Range("T2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_1885/kartor/A 34293-2023 karta.png", "A 34293-2023")'
$ws.Range("V2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_1885/klagomål/A 34293-2023 fsc-klagomål.docx", "A 34293-2023")'
$ws.Range("W2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_1885/klagomålsmail/A 34293-2023 fsc-klagomål mail.docx", "A 34293-2023")'
$ws.Range("X2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_1885/tillsyn/A 34293-2023 tillsynsbegäran.docx", "A 34293-2023")'
$ws.Range("Y2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_1885/ti,llsynsmail/A 34293-2023 tillsynsbegäran mail.docx", "A 34293-2023")'
